# Insert a new data row at row 446 (pushes existing rows 446-554 down to 447-555)
# and populate it with the new weekly price record, matching the commit message
# "Fruta / hortaliza, semanal" (weekly fruit/vegetable price update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("446:446").Insert()

$ws.Range("A446").Value = 7
$ws.Range("B446").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C446").Value = "Ñuble"
$ws.Range("D446").Value = 45244
$ws.Range("E446").Value = 16
$ws.Range("F446").Value = 100112003
$ws.Range("G446").Value = "Ajo"
$ws.Range("H446").Value = "Chino"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 100
$ws.Range("K446").Value = 24000
$ws.Range("L446").Value = 24000
$ws.Range("M446").Value = 24000
$ws.Range("N446").Value = "`$/caja 10 kilos"
$ws.Range("O446").Value = "China"
$ws.Range("P446").Value = 2400
$ws.Range("Q446").Value = 10
$ws.Range("R446").Value = "Hortaliza"

$ws.Range("D446").NumberFormat = "YYYY-MM-DD HH:MM:SS"
